$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (changed) date column C for rows 2-28: 45475 -> 45476
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45476
}

# 2) Row 28 gains an explicit row height (15, custom) - it previously had none
$ws.Rows.Item(28).RowHeight = 15

# 3) Append a new row 29 with the new record
$ws.Range("A29").Value = "A 27830-2024"
$ws.Range("B29").Value = 45475
$ws.Range("C29").Value = 45476
$ws.Range("B29:C29").NumberFormat = "YYYY-MM-DD"
$ws.Range("D29").Value = "OKÄNT"
$ws.Range("E29").Value = "OKÄNT"
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("R29").WrapText = $true
